# Project Timeline.xlsx — "Updated Documentation,: Financial plan, House of
# Quality, Engineering Requirements"
#
# Functional changes applied:
#  1. Make "Fire Alarm" (sheet 1) the active/selected sheet instead of
#     "Teaching Assistant" (sheet 2).
#  2. Fire Alarm page setup: turn on "Fit to page" and set print scale to 74%
#     (keep landscape orientation).
#  3. Teaching Assistant page setup: switch orientation from portrait to
#     landscape, turn on "Fit to page" and set print scale to 71%.

$wb = $excel.ActiveWorkbook

$fireAlarm = $wb.Worksheets.Item("Fire Alarm")
$teachingAssistant = $wb.Worksheets.Item("Teaching Assistant")

# --- Page setup: Fire Alarm (sheet1.xml) ---------------------------------
$fireAlarm.PageSetup.Zoom = 74
$fireAlarm.PageSetup.FitToPagesWide = $False

# --- Page setup: Teaching Assistant (sheet2.xml) --------------------------
$teachingAssistant.PageSetup.Orientation = 2
$teachingAssistant.PageSetup.Zoom = 71
$teachingAssistant.PageSetup.FitToPagesWide = $False

# --- Active sheet: Fire Alarm becomes the selected tab ---------------------
$fireAlarm.Activate()
